$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): Wins, Losses, Ties in AC1, AD1, AE1
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header style (bold, centered, bordered) from an existing header cell (AB1)
$ws.Range("AB1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Re-set values after paste (paste formats shouldn't touch values, but just to be safe)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Data rows 2-48: Wins=67, Losses=95, Ties=0
for ($r = 2; $r -le 48; $r++) {
    $ws.Range("AC$r").Value = 67
    $ws.Range("AD$r").Value = 95
    $ws.Range("AE$r").Value = 0
}

$excel.DisplayAlerts = $false
$wb.Save()
